$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.09250000000001
$ws.Range("A9").Value = -21.9453
$ws.Range("E12").Value = 18.17850000000003
$ws.Range("A13").Value = -22.1329
$ws.Range("E14").Value = 16.84940000000001
$ws.Range("A16").Value = -21.61319999999999
$ws.Range("A18").Value = -22.16740000000001
$ws.Range("E19").Value = 16.6846
$ws.Range("A20").Value = -21.16019999999997
$ws.Range("A26").Value = -21.06639999999996
$ws.Range("E26").Value = 15.99799999999999
$ws.Range("A27").Value = -22.02089999999998
$ws.Range("E27").Value = 16.56389999999999
$ws.Range("A29").Value = -20.65829999999997
$ws.Range("E29").Value = 16.87030000000001
$ws.Range("A35").Value = -21.22259999999998
$ws.Range("A36").Value = -21.35909999999999
$ws.Range("E37").Value = 16.72390000000001
$ws.Range("E38").Value = 16.3339
$ws.Range("A45").Value = -21.49919999999999
$ws.Range("E47").Value = 16.5385
$ws.Range("E51").Value = 17.1923
$ws.Range("E52").Value = 17.18000000000001
$ws.Range("A55").Value = -22.1202
$ws.Range("E55").Value = 16.47250000000001
$ws.Range("A57").Value = -22.215
$ws.Range("A69").Value = -21.68299999999996
$ws.Range("E69").Value = 17.22360000000002
$ws.Range("E70").Value = 18.08540000000002
$ws.Range("A76").Value = -19.58519999999999
$ws.Range("E76").Value = 16.45699999999999
$ws.Range("A78").Value = -19.89479999999998
$ws.Range("E81").Value = 16.56920000000001
$ws.Range("A82").Value = -21.98750000000001
$ws.Range("A83").Value = -21.9311
$ws.Range("E83").Value = 16.62299999999999
$ws.Range("A93").Value = -20.67019999999998
$ws.Range("E94").Value = 18.72570000000002
$ws.Range("A97").Value = -21.82010000000001
$ws.Range("E100").Value = 16.3321
$ws.Range("E102").Value = 16.7873
